$d = $word.ActiveDocument

# Locate the very last paragraph in the document body ("Move enemy/size
# database into a file" + trailing tab) and append three new list items
# right after it, before the section properties.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)

# Create a fresh paragraph to host the inserted XML (InsertXML replaces
# the contents of the exact range it is called on, so we give it an
# empty paragraph of its own rather than risk clobbering the preceding
# "Move enemy/size database into a file" paragraph).
$null = $r.InsertParagraphAfter()
$target = $d.Paragraphs.Last.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Adobe Caslon Pro" w:hAnsi="Adobe Caslon Pro"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Adobe Caslon Pro" w:hAnsi="Adobe Caslon Pro"/></w:rPr><w:t xml:space="preserve">Have a crack at using some OpenGL </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Adobe Caslon Pro" w:hAnsi="Adobe Caslon Pro"/></w:rPr><w:t>shaders</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Adobe Caslon Pro" w:hAnsi="Adobe Caslon Pro"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Adobe Caslon Pro" w:hAnsi="Adobe Caslon Pro"/></w:rPr><w:t>Distortion exploding effect</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Adobe Caslon Pro" w:hAnsi="Adobe Caslon Pro"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Adobe Caslon Pro" w:hAnsi="Adobe Caslon Pro"/></w:rPr><w:t>Instance particles?</w:t></w:r></w:p>
'@

$null = $target.InsertXML($xml)
